# Apply edits described by the diff: zero out several parameter cells on
# the "agent1" sheet (rows 42-44, columns L/M/P/Q/R). Because "agent2" and
# "agent3" hold formulas like "=agent1!Q42" these changes ripple through on
# recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("agent1")

$ws.Range("Q42").Value = 0
$ws.Range("R42").Value = 0

$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0
$ws.Range("P43").Value = 0
$ws.Range("Q43").Value = 0
$ws.Range("R43").Value = 0

$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("P44").Value = 0
$ws.Range("Q44").Value = 0
$ws.Range("R44").Value = 0
